$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valkyrie Stage-I")
Write-Host $ws.Hyperlinks.Count
foreach ($hl in $ws.Hyperlinks) {
    Write-Host $hl.Range.Address() "|" $hl.TextToDisplay "|" $hl.Address
}
